$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.125.30"
$ws.Range("E2").Value = "  +6.66%  "

$ws.Range("D3").Value = "3.520.03"
$ws.Range("E3").Value = "  +9.78%  "

$ws.Range("E4").Value = "  +0.28%  "

$ws.Range("D5").Value = "'189.55"
$ws.Range("E5").Value = "  +9.21%  "

$ws.Range("D6").Value = "'552.41"
$ws.Range("E6").Value = "  +3.87%  "

$ws.Range("D7").Value = "3.521.41"
$ws.Range("E7").Value = "  +9.75%  "

$ws.Range("D8").Value = "'0.606"
$ws.Range("E8").Value = "  +1.91%  "

$ws.Range("E9").Value = "  +0.04%  "

$ws.Range("D10").Value = "'0.633"
$ws.Range("E10").Value = "  +3.92%  "

$ws.Range("E11").Value = "  +15.24%  "

$ws.Range("D12").Value = "'54.95"
$ws.Range("E12").Value = "  +2.34%  "

$ws.Range("E13").Value = "  +6.43%  "

$ws.Range("D14").Value = "'9.36"
$ws.Range("E14").Value = "  +2.80%  "

$ws.Range("D15").Value = "4.088.96"
$ws.Range("E15").Value = "  +10.08%  "

$ws.Range("D16").Value = "3.527.92"
$ws.Range("E16").Value = "  +10.78%  "

$ws.Range("E17").Value = "  +3.26%  "

$ws.Range("D18").Value = "67.237.02"
$ws.Range("E18").Value = "  +7.37%  "

$ws.Range("D19").Value = "'18.22"
$ws.Range("E19").Value = "  +5.09%  "

$ws.Range("D20").Value = "'11.90"
$ws.Range("E20").Value = "  +7.11%  "

$ws.Range("D21").Value = "'0.996"
$ws.Range("E21").Value = "  +2.69%  "

$ws.Range("D22").Value = "'430.32"
$ws.Range("E22").Value = "  +17.08%  "

$ws.Range("E23").Value = "  +3.47%  "

$ws.Range("D24").Value = "'85.01"
$ws.Range("E24").Value = "  +4.46%  "

$ws.Range("D25").Value = "'4.17"
$ws.Range("E25").Value = "  +7.44%  "

$ws.Range("D26").Value = "'11.18"
$ws.Range("E26").Value = "  -0.93%  "

$ws.Range("E27").Value = "  +9.02%  "

$ws.Range("E28").Value = "  +5.43%  "

$ws.Range("E29").Value = "  +9.33%  "

$ws.Range("D30").Value = "'30.24"
$ws.Range("E30").Value = "  +6.17%  "

$ws.Range("D31").Value = "'648.93"
$ws.Range("E31").Value = "  +0.46%  "

$ws.Range("D32").Value = "'6.69"
$ws.Range("E32").Value = "  +2.23%  "

$ws.Range("D33").Value = "'11.72"
$ws.Range("E33").Value = "  +3.18%  "

$ws.Range("E34").Value = "  +4.60%  "

$ws.Range("D35").Value = "'59.30"
$ws.Range("E35").Value = "  +4.15%  "

$ws.Range("B36").Value = "InjectiveProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D36").Value = "'38.60"
$ws.Range("E36").Value = "  +4.00%  "

$ws.Range("B37").Value = "PEPE"
$ws.Range("C37").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D37").Value = "0.0₃0817"
$ws.Range("E37").Value = "  +14.08%  "

$ws.Range("E38").Value = "  -0.30%  "

$ws.Range("D39").Value = "'0.391"
$ws.Range("E39").Value = "  +3.37%  "

$ws.Range("E40").Value = "  +14.75%  "

$ws.Range("E41").Value = "  +14.15%  "

$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  +0.55%  "

$ws.Range("D43").Value = "3.024.46"
$ws.Range("E43").Value = "  +4.90%  "

$ws.Range("D44").Value = "'2.65"
$ws.Range("E44").Value = "  +3.93%  "

$ws.Range("E45").Value = "  +9.15%  "

$ws.Range("E46").Value = "  +10.16%  "

$ws.Range("D47").Value = "'3.32"
$ws.Range("E47").Value = "  +11.87%  "

$ws.Range("E48").Value = "  +6.02%  "

$ws.Range("D49").Value = "'0.130"
$ws.Range("E49").Value = "  +5.23%  "

$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").Value = "'142.88"
$ws.Range("E50").Value = "  +6.41%  "

$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").Value = "'8.71"
$ws.Range("E51").Value = "  +12.88%  "
